$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46022
$ws.Range("B2").Value = 10918.5663016096
$ws.Range("C2").Value = 11395.2765385751
$ws.Range("D2").Value = 19632.26
$ws.Range("E2").Value = 8155.15888828042
$ws.Range("F2").Value = -3.40935721435365

$ws.Range("A3").Value = 46023
$ws.Range("B3").Value = 5104.69244805461
$ws.Range("C3").Value = 9146.37939665343
$ws.Range("D3").Value = 12075.86
$ws.Range("E3").Value = 8708.69966055211
$ws.Range("F3").Value = 240.800794050231

$ws.Range("A4").Value = 46024
$ws.Range("B4").Value = 11522.5601961688
$ws.Range("C4").Value = 10972.8086095921
$ws.Range("D4").Value = 12075.86
$ws.Range("E4").Value = 8923.262080073
$ws.Range("F4").Value = 325.842112069379

$ws.Range("A5").Value = 46025
$ws.Range("B5").Value = 4951.25180548264
$ws.Range("C5").Value = 7902.66338132959
$ws.Range("D5").Value = 12075.86
$ws.Range("E5").Value = 8230.74287508738
$ws.Range("F5").Value = 169.064427350707

$ws.Range("A6").Value = 46026
$ws.Range("B6").Value = 4869.75255135434
$ws.Range("C6").Value = 8089.3265103794
$ws.Range("D6").Value = 12075.86
$ws.Range("E6").Value = 8269.79579047536
$ws.Range("F6").Value = 178.469262535615

$ws.Range("A7").Value = 46027
$ws.Range("B7").Value = 13228.8169926471
$ws.Range("C7").Value = 13059.9218961483
$ws.Range("D7").Value = 12075.86
$ws.Range("E7").Value = 9108.31864906552
$ws.Range("F7").Value = 420.515856050577

$ws.Range("A8").Value = 46028
$ws.Range("B8").Value = 5295.53439346162
$ws.Range("C8").Value = 9514.77467648482
$ws.Range("D8").Value = 12075.86
$ws.Range("E8").Value = 8688.20133476951
$ws.Range("F8").Value = 255.29650046893

$ws.Range("A9").Value = 46029
$ws.Range("B9").Value = 13228.8169926471
$ws.Range("C9").Value = 13574.2067762396
$ws.Range("D9").Value = 12075.86
$ws.Range("E9").Value = 9108.31864906552
$ws.Range("F9").Value = 441.944392721046

$ws.Range("A10").Value = 46030
$ws.Range("B10").Value = 13228.8169926471
$ws.Range("C10").Value = 13425.7156596898
$ws.Range("D10").Value = 12075.86
$ws.Range("E10").Value = 9108.31864906552
$ws.Range("F10").Value = 435.757262864806

$ws.Range("A11").Value = 46031
$ws.Range("B11").Value = 13228.8169926471
$ws.Range("C11").Value = 12294.424660546
$ws.Range("D11").Value = 12075.86
$ws.Range("E11").Value = 9108.31864906552
$ws.Range("F11").Value = 388.620137900481

$ws.Range("A12").Value = 46032
$ws.Range("B12").Value = 5461.34628757431
$ws.Range("C12").Value = 8574.93108425919
$ws.Range("D12").Value = 12075.86
$ws.Range("E12").Value = 8700.47673030118
$ws.Range("F12").Value = 216.647825606682

$ws.Range("A13").Value = 46033
$ws.Range("B13").Value = 5295.53439346162
$ws.Range("C13").Value = 8528.09924486521
$ws.Range("D13").Value = 12075.86
$ws.Range("E13").Value = 8688.20133476951
$ws.Range("F13").Value = 214.185024151447

$ws.Range("A14").Value = 46034
$ws.Range("B14").Value = 13381.7270213684
$ws.Range("C14").Value = 12813.0904514226
$ws.Range("D14").Value = 12075.86
$ws.Range("E14").Value = 9107.03037487135
$ws.Range("F14").Value = 410.177534428916

$ws.Range("A15").Value = 46035
$ws.Range("B15").Value = 13381.7270213684
$ws.Range("C15").Value = 13537.3175641503
$ws.Range("D15").Value = 12075.86
$ws.Range("E15").Value = 9107.03037487135
$ws.Range("F15").Value = 440.353664125902

